# Update "want to go" counts (and min-price for a couple rows) on the
# "展览" and "全部类型" sheets — both carry the same underlying rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        # 合肥·第九届环形宇宙动漫游戏嘉年华 / 合肥·MAX特摄同人only2.0 live on rows 11/12 here
        $rumorRow1 = 11
        $rumorRow2 = 12
    } else {
        # same two rows live on rows 14/15 on the combined "全部类型" sheet
        $rumorRow1 = 14
        $rumorRow2 = 15
    }

    $ws.Range("F2").Value = 302
    $ws.Range("F4").Value = 8052
    $ws.Range("F5").Value = 5870
    $ws.Range("F7").Value = 89

    $ws.Cells.Item($rumorRow1, 6).Value = 400
    $ws.Cells.Item($rumorRow1, 7).Value = 72
    $ws.Cells.Item($rumorRow2, 6).Value = 67
}
